$wb = $excel.ActiveWorkbook

# Locate the existing "ShuntSw" sheet; the new "Shunt" sheet is derived from
# it (same column layout minus the switching-specific columns), so copy it
# and place the copy immediately before "ShuntSw" -- this also reproduces
# the sheetId/rId shift seen in the target workbook.
$shuntSw = $wb.Worksheets.Item("ShuntSw")
$shuntSw.Copy($shuntSw)

$shunt = $wb.Worksheets.Item(6)
$shunt.Name = "Shunt"

# The ShuntSw-only columns (gs, bs, ns, vref, dv, dt) don't exist on the
# Shunt model; drop them so the sheet ends at column J.
$shunt.Columns("K:P").Delete()

# Rename the two device rows from ShuntSw_n to Shunt_n ...
$shunt.Range("B2").Value = "Shunt_1"
$shunt.Range("D2").Value = "Shunt_1"
$shunt.Range("B3").Value = "Shunt_2"
$shunt.Range("D3").Value = "Shunt_2"

# ... the leftover bold "Normal 2" formatting on a couple of the copied
# name cells doesn't apply to the new data, so strip it back to default.
$shunt.Range("D2").ClearFormats()
$shunt.Range("B3").ClearFormats()
$shunt.Range("D3").ClearFormats()

# Shunt_1's rated power differs from the ShuntSw device it was copied from.
$shunt.Range("F2").Value = 100

# Restore the plain selection the new sheet should carry.
$shunt.Range("I13").Select()
